# Finalized Experiments with Participant Generation
# Renames the task-order sheets to the new generated task-order ids, and
# replaces the generated stimulus-file names in each sheet accordingly.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (order matches sheetId 1..5 / rId1..rId5) ---
$wsGNG  = $wb.Worksheets.Item(1)
$wsNB   = $wb.Worksheets.Item(2)
$wsRS   = $wb.Worksheets.Item(3)
$wsTOL  = $wb.Worksheets.Item(4)
$wsvSAT = $wb.Worksheets.Item(5)

$wsGNG.Name  = "GNG_TO-1650291249263069"
$wsNB.Name   = "NB_TO-16502912516043196"
$wsRS.Name   = "RS_TO-16502912516073222"
$wsTOL.Name  = "TOL_TO-16502912516823235"
$wsvSAT.Name = "vSAT_TO-16502912517773306"

# --- Sheet 1 (GNG): update generated stim file names ---
$wsGNG.Range("B2").Value = "go_stims-1650291249207066.csv"
$wsGNG.Range("B3").Value = "GNG_stims-16502912492330682.csv"
$wsGNG.Range("B4").Value = "go_stims-16502912492340732.csv"
$wsGNG.Range("B5").Value = "GNG_stims-16502912492620704.csv"

# --- Sheet 2 (NB): update generated stim file names ---
$wsNB.Range("B2").Value = "ZB-match_4-16502912496149418.csv"
$wsNB.Range("B3").Value = "TB-1650291251335322.csv"
$wsNB.Range("B4").Value = "OB-165029125122132.csv"
$wsNB.Range("B5").Value = "OB-16502912499540627.csv"
$wsNB.Range("B6").Value = "ZB-match_6-16502912496529353.csv"
$wsNB.Range("B7").Value = "TB-16502912515783215.csv"
$wsNB.Range("B8").Value = "TB-16502912514393523.csv"
$wsNB.Range("B9").Value = "ZB-match_1-16502912496869433.csv"
$wsNB.Range("B10").Value = "OB-16502912508183222.csv"

# --- Sheet 3 (RS): no content changes ---

# --- Sheet 4 (TOL): update generated stim file names ---
$wsTOL.Range("B2").Value = "MM_stims-1650291251636384.csv"
$wsTOL.Range("B3").Value = "ZM_stims-16502912516103272.csv"
$wsTOL.Range("B4").Value = "MM_stims-16502912516663554.csv"
$wsTOL.Range("B5").Value = "ZM_stims-16502912516383367.csv"
$wsTOL.Range("B6").Value = "MM_stims-16502912516813226.csv"
$wsTOL.Range("B7").Value = "ZM_stims-16502912516673226.csv"

# --- Sheet 5 (vSAT): update generated stim file names ---
$wsvSAT.Range("B2").Value = "vSAT_stims-16502912517293212.csv"
$wsvSAT.Range("B3").Value = "vSAT_stims-16502912517623222.csv"
$wsvSAT.Range("B4").Value = "SAT_stims-16502912517133224.csv"
$wsvSAT.Range("B5").Value = "SAT_stims-16502912516883223.csv"
